$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 623.6875
$ws.Range("I103").Value = 656.6
$ws.Range("J103").Value = 608.7273
$ws.Range("K103").Value = 1969.8
$ws.Range("L103").Value = 1826.1819
$ws.Range("M103").Value = -1383.8
$ws.Range("N103").Value = -2998.1819
$ws.Range("H113").Value = 3921.6191
$ws.Range("I113").Value = 5181.625
$ws.Range("J113").Value = 3146.2307
$ws.Range("K113").Value = 5181.625
$ws.Range("L113").Value = 3146.2307
$ws.Range("M113").Value = -1927.625
$ws.Range("N113").Value = -9654.2307
$ws.Range("H129").Value = 25642088
$ws.Range("J129").Value = 1026.25
$ws.Range("L129").Value = 3078.75
$ws.Range("N129").Value = -13078.75
$ws.Range("H132").Value = 8776395
$ws.Range("I132").Value = 9438651
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 28315953
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -28313423
$ws.Range("N132").Value = -9560
$ws.Range("H135").Value = 873.75
$ws.Range("I135").Value = 873.75
$ws.Range("K135").Value = 7863.75
$ws.Range("M135").Value = -5328.75
$ws.Range("H138").Value = 3215.6567
$ws.Range("I138").Value = 1688.8518
$ws.Range("J138").Value = 4246.25
$ws.Range("K138").Value = 5066.555399999999
$ws.Range("L138").Value = 12738.75
$ws.Range("M138").Value = 73.44460000000072
$ws.Range("N138").Value = -23018.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9294.583000000001
$ws.Range("I32").Value = 9066.173000000001
$ws.Range("J32").Value = 10528
$ws.Range("K32").Value = 9066.173000000001
$ws.Range("L32").Value = 10528
$ws.Range("M32").Value = -8779.173000000001
$ws.Range("N32").Value = -11102
$ws.Range("H45").Value = 10449212
$ws.Range("I45").Value = 12860168
$ws.Range("J45").Value = 1733
$ws.Range("K45").Value = 12860168
$ws.Range("L45").Value = 1733
$ws.Range("M45").Value = -12859791
$ws.Range("N45").Value = -2487
$ws.Range("H61").Value = 1737.2046
$ws.Range("I61").Value = 1700.3462
$ws.Range("J61").Value = 1790.4445
$ws.Range("K61").Value = 1700.3462
$ws.Range("L61").Value = 1790.4445
$ws.Range("M61").Value = -1488.3462
$ws.Range("N61").Value = -2214.4445
$ws.Range("H63").Value = 1699.4286
$ws.Range("I63").Value = 1598.9231
$ws.Range("J63").Value = 3006
$ws.Range("K63").Value = 1598.9231
$ws.Range("L63").Value = 3006
$ws.Range("M63").Value = -912.9231
$ws.Range("N63").Value = -4378
$ws.Range("H66").Value = 1699.4286
$ws.Range("I66").Value = 1598.9231
$ws.Range("J66").Value = 3006
$ws.Range("K66").Value = 7994.6155
$ws.Range("L66").Value = 15030
$ws.Range("M66").Value = -4562.6155
$ws.Range("N66").Value = -21894
$ws.Range("H102").Value = 1585
$ws.Range("I102").Value = 1585
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1585
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 37
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 2761.1875
$ws.Range("I122").Value = 2877.6667
$ws.Range("J122").Value = 1014
$ws.Range("K122").Value = 8633.000100000001
$ws.Range("L122").Value = 3042
$ws.Range("M122").Value = -6183.000100000001
$ws.Range("N122").Value = -7942
$ws.Range("H132").Value = 3886.4102
$ws.Range("I132").Value = 978.8125
$ws.Range("J132").Value = 17178.285
$ws.Range("K132").Value = 2936.4375
$ws.Range("L132").Value = 51534.855
$ws.Range("M132").Value = -406.4375
$ws.Range("N132").Value = -56594.855
$ws.Range("H136").Value = 1737.2046
$ws.Range("I136").Value = 1700.3462
$ws.Range("J136").Value = 1790.4445
$ws.Range("K136").Value = 5101.0386
$ws.Range("L136").Value = 5371.333500000001
$ws.Range("M136").Value = -2551.0386
$ws.Range("N136").Value = -10471.3335

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3575.5881
$ws.Range("I20").Value = 3337.8
$ws.Range("J20").Value = 3915.2856
$ws.Range("K20").Value = 3337.8
$ws.Range("L20").Value = 3915.2856
$ws.Range("M20").Value = -3090.8
$ws.Range("N20").Value = -4409.2856
$ws.Range("H107").Value = 23810498
$ws.Range("I107").Value = 31250718
$ws.Range("K107").Value = 31250718
$ws.Range("M107").Value = -31248798
$ws.Range("H122").Value = 40494.43
$ws.Range("J122").Value = 40494.43
$ws.Range("L122").Value = 40494.43
$ws.Range("N122").Value = -50294.43
$ws.Range("H134").Value = 2590729.8
$ws.Range("I134").Value = 6697.227
$ws.Range("J134").Value = 5297811.5
$ws.Range("K134").Value = 20091.681
$ws.Range("L134").Value = 15893434.5
$ws.Range("M134").Value = -17556.681
$ws.Range("N134").Value = -15898504.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1786.9688
$ws.Range("I134").Value = 1929
$ws.Range("K134").Value = 5787
$ws.Range("M134").Value = -3252

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 13095760
$ws.Range("I113").Value = 8333834
$ws.Range("J113").Value = 25000574
$ws.Range("K113").Value = 25001502
$ws.Range("L113").Value = 75001722
$ws.Range("M113").Value = -24999332
$ws.Range("N113").Value = -75006062
$ws.Range("H131").Value = 708.3200000000001
$ws.Range("J131").Value = 734.8602
$ws.Range("L131").Value = 2204.5806
$ws.Range("N131").Value = -12284.5806

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1170.6471
$ws.Range("I97").Value = 1102.7273
$ws.Range("J97").Value = 1295.1666
$ws.Range("K97").Value = 1102.7273
$ws.Range("L97").Value = 1295.1666
$ws.Range("M97").Value = -606.7273
$ws.Range("N97").Value = -2287.1666
$ws.Range("H132").Value = 3917.9285
$ws.Range("I132").Value = 991.5517
$ws.Range("K132").Value = 2974.6551
$ws.Range("M132").Value = -444.6550999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6258.024
$ws.Range("I132").Value = 2379.6667
$ws.Range("J132").Value = 9166.791999999999
$ws.Range("K132").Value = 7139.000100000001
$ws.Range("L132").Value = 27500.376
$ws.Range("M132").Value = -4609.000100000001
$ws.Range("N132").Value = -32560.376
$ws.Range("H136").Value = 2843.0625
$ws.Range("I136").Value = 2312.2896
$ws.Range("J136").Value = 4860
$ws.Range("K136").Value = 6936.8688
$ws.Range("L136").Value = 14580
$ws.Range("M136").Value = -4386.8688
$ws.Range("N136").Value = -19680

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 27778366
$ws.Range("I81").Value = 31250524
$ws.Range("J81").Value = 1101
$ws.Range("K81").Value = 62501048
$ws.Range("L81").Value = 2202
$ws.Range("M81").Value = -62499987
$ws.Range("N81").Value = -4324
$ws.Range("H84").Value = 27778366
$ws.Range("I84").Value = 31250524
$ws.Range("J84").Value = 1101
$ws.Range("K84").Value = 312505240
$ws.Range("L84").Value = 11010
$ws.Range("M84").Value = -312499936
$ws.Range("N84").Value = -21618
$ws.Range("H132").Value = 13174130
$ws.Range("I132").Value = 19251718
$ws.Range("J132").Value = 6024.2915
$ws.Range("K132").Value = 57755154
$ws.Range("L132").Value = 18072.8745
$ws.Range("M132").Value = -57752624
$ws.Range("N132").Value = -23132.8745
$ws.Range("H136").Value = 2903.8667
$ws.Range("I136").Value = 4427.4
$ws.Range("J136").Value = 1380.3334
$ws.Range("K136").Value = 13282.2
$ws.Range("L136").Value = 4141.0002
$ws.Range("M136").Value = -10732.2
$ws.Range("N136").Value = -9241.0002
